$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "43.718.81"
$ws.Range("E2").Value = "  -0.55%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.295.42"
$ws.Range("E3").Value = "  +1.77%  "
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "231.77"
$ws.Range("E5").Value = "  +1.00%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.643"
$ws.Range("E6").Value = "  +1.62%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "64.41"
$ws.Range("E7").Value = "  +1.53%  "
$ws.Range("E8").Value = "  -0.18%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.440"
$ws.Range("E9").Value = "  -0.18%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0964"
$ws.Range("E10").Value = "  -5.05%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "56.88"
$ws.Range("E11").Value = "  -0.38%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "26.60"
$ws.Range("E12").Value = "  +2.28%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "2.628.33"
$ws.Range("E14").Value = "  +1.40%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "15.17"
$ws.Range("E15").Value = "  -2.77%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "6.09"
$ws.Range("E16").Value = "  +0.04%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.831"
$ws.Range("E17").Value = "  -1.61%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.287.17"
$ws.Range("E18").Value = "  +1.13%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "43.602.83"
$ws.Range("E19").Value = "  -0.76%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.0₃0973"
$ws.Range("E20").Value = "  -3.16%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "73.28"
$ws.Range("E21").Value = "  +0.22%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.14"
$ws.Range("E22").Value = "  +2.21%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "248.26"
$ws.Range("E23").Value = "  -1.38%  "
$ws.Range("E24").Value = "  +0.08%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "3.72"
$ws.Range("E25").Value = "  +11.68%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.45"
$ws.Range("E26").Value = "  +1.38%  "
$ws.Range("E27").Value = "  -1.90%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "9.75"
$ws.Range("E28").Value = "  -2.13%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "174.10"
$ws.Range("E29").Value = "  +1.25%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "21.97"
$ws.Range("E30").Value = "  +6.20%  "
$ws.Range("B31").Value = "Kaspa"
$ws.Range("C31").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.131"
$ws.Range("E31").Value = "  -4.18%  "
$ws.Range("B32").Value = "ImmutableX"
$ws.Range("C32").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.42"
$ws.Range("E32").Value = "  +3.38%  "
$ws.Range("E33").Value = "  +1.28%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.95"
$ws.Range("E34").Value = "  +5.27%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0682"
$ws.Range("E35").Value = "  -0.26%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "4.95"
$ws.Range("E36").Value = "  +2.74%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "6.49"
$ws.Range("E37").Value = "  -0.41%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.60"
$ws.Range("E38").Value = "  -5.12%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.31"
$ws.Range("E39").Value = "  +0.70%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0251"
$ws.Range("E40").Value = "  -2.02%  "
$ws.Range("E41").Value = "  +0.08%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "8.78"
$ws.Range("E42").Value = "  +6.89%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "17.18"
$ws.Range("E43").Value = "  -0.81%  "
$ws.Range("E44").Value = "  +2.42%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "97.31"
$ws.Range("E45").Value = "  +0.04%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "10.20"
$ws.Range("E46").Value = "  +7.33%  "
$ws.Range("B47").Value = "Cronos"
$ws.Range("C47").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0944"
$ws.Range("E47").Value = "  -2.36%  "
$ws.Range("B48").Value = "TrustWalletToken"
$ws.Range("C48").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.18"
$ws.Range("E48").Value = "  -0.36%  "
$ws.Range("E49").Value = "  -3.86%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.432.84"
$ws.Range("E50").Value = "  -0.25%  "
$ws.Range("B51").Value = "ARBITRUM"
$ws.Range("C51").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.09"
$ws.Range("E51").Value = "  +1.43%  "
